# Scheduled market-data refresh: update cached average-price / Leve-profit
# figures (columns H-N) on a handful of rows across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets to the latest pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 11226.286
$ws.Range("I132").Value = 9288.091
$ws.Range("J132").Value = 18333
$ws.Range("K132").Value = 27864.273
$ws.Range("L132").Value = 54999
$ws.Range("M132").Value = -25334.273
$ws.Range("N132").Value = -60059

# Row 135
$ws.Range("H135").Value = 4941.9653
$ws.Range("I135").Value = 1827.6
$ws.Range("J135").Value = 11862.777
$ws.Range("K135").Value = 16448.4
$ws.Range("L135").Value = 106764.993
$ws.Range("M135").Value = -13913.4

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6398.483
$ws.Range("I61").Value = 7236
$ws.Range("J61").Value = 5718
$ws.Range("K61").Value = 7236
$ws.Range("L61").Value = 5718
$ws.Range("M61").Value = -7024

# Row 122
$ws.Range("H122").Value = 6770.067
$ws.Range("I122").Value = 6242.1665
$ws.Range("J122").Value = 7122
$ws.Range("K122").Value = 18726.4995
$ws.Range("L122").Value = 21366
$ws.Range("M122").Value = -16276.4995

# Row 136
$ws.Range("H136").Value = 6398.483
$ws.Range("I136").Value = 7236
$ws.Range("J136").Value = 5718
$ws.Range("K136").Value = 21708
$ws.Range("L136").Value = 17154
$ws.Range("M136").Value = -19158

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 792.61536
$ws.Range("I20").Value = 792.61536
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 792.61536
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -545.61536
$ws.Range("N20").Value = $null

# Row 99
$ws.Range("H99").Value = 4419.4
$ws.Range("I99").Value = 3699
$ws.Range("J99").Value = 4899.6665
$ws.Range("K99").Value = 3699
$ws.Range("L99").Value = 4899.6665
$ws.Range("M99").Value = -2201
$ws.Range("N99").Value = -7895.6665

# Row 105
$ws.Range("H105").Value = 38465260

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 14348.3
$ws.Range("I62").Value = 16915.5
$ws.Range("J62").Value = 10497.5
$ws.Range("K62").Value = 16915.5
$ws.Range("L62").Value = 10497.5
$ws.Range("M62").Value = -16291.5
$ws.Range("N62").Value = -11745.5

# Row 65
$ws.Range("H65").Value = 14348.3
$ws.Range("I65").Value = 16915.5
$ws.Range("J65").Value = 10497.5
$ws.Range("K65").Value = 84577.5
$ws.Range("L65").Value = 52487.5
$ws.Range("M65").Value = -81457.5
$ws.Range("N65").Value = -58727.5

# Row 132
$ws.Range("H132").Value = 12317.23
$ws.Range("I132").Value = 17483.25
$ws.Range("J132").Value = 10021.223
$ws.Range("K132").Value = 52449.75
$ws.Range("L132").Value = 30063.669
$ws.Range("M132").Value = -49919.75
$ws.Range("N132").Value = -35123.669

# Row 134
$ws.Range("H134").Value = 24500.467
$ws.Range("I134").Value = 20498.5
$ws.Range("J134").Value = 25955.727
$ws.Range("K134").Value = 61495.5
$ws.Range("L134").Value = 77867.181
$ws.Range("M134").Value = -58960.5
$ws.Range("N134").Value = -82937.181

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Range("H62").Value = 15324.875
$ws.Range("I62").Value = 5999.5
$ws.Range("J62").Value = 18433.334
$ws.Range("K62").Value = 17998.5
$ws.Range("L62").Value = 55300.00199999999
$ws.Range("M62").Value = -17312.5
$ws.Range("N62").Value = -56672.00199999999

# Row 65
$ws.Range("H65").Value = 15324.875
$ws.Range("I65").Value = 5999.5
$ws.Range("J65").Value = 18433.334
$ws.Range("K65").Value = 53995.5
$ws.Range("L65").Value = 165900.006
$ws.Range("M65").Value = -50563.5
$ws.Range("N65").Value = -172764.006

# Row 98
$ws.Range("H98").Value = 3831.6667
$ws.Range("I98").Value = 3501.5
$ws.Range("J98").Value = 3996.75
$ws.Range("K98").Value = 10504.5
$ws.Range("L98").Value = 11990.25
$ws.Range("M98").Value = -9006.5
$ws.Range("N98").Value = -14986.25

# Row 131
$ws.Range("H131").Value = 37683556
$ws.Range("I131").Value = 59259680
$ws.Range("J131").Value = 23813188
$ws.Range("K131").Value = 177779040
$ws.Range("L131").Value = 71439564
$ws.Range("M131").Value = -177774000
$ws.Range("N131").Value = -71449644

# Row 136
$ws.Range("H136").Value = 8774417
$ws.Range("I136").Value = 12822409
$ws.Range("J136").Value = 3766.3333
$ws.Range("K136").Value = 38467227
$ws.Range("L136").Value = 11298.9999
$ws.Range("M136").Value = -38462127
$ws.Range("N136").Value = -21498.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 59875
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 59875
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 59875
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -61871

# Row 83
$ws.Range("H83").Value = 59875
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 59875
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 299375
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -309359

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4682.171
$ws.Range("I7").Value = 4731.963
$ws.Range("J7").Value = 4586.143
$ws.Range("K7").Value = 4731.963
$ws.Range("L7").Value = 4586.143
$ws.Range("M7").Value = -4619.963

# Row 22
$ws.Range("H22").Value = 1250
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -2590

# Row 27
$ws.Range("H27").Value = 1250
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -2214

# Row 68
$ws.Range("H68").Value = 2294
$ws.Range("I68").Value = 2538.3845
$ws.Range("J68").Value = 1499.75
$ws.Range("K68").Value = 2538.3845
$ws.Range("L68").Value = 1499.75
$ws.Range("M68").Value = -1789.3845
$ws.Range("N68").Value = -2997.75

# Row 71
$ws.Range("H71").Value = 2294
$ws.Range("I71").Value = 2538.3845
$ws.Range("J71").Value = 1499.75
$ws.Range("K71").Value = 12691.9225
$ws.Range("L71").Value = 7498.75
$ws.Range("M71").Value = -8947.922500000001
$ws.Range("N71").Value = -14986.75

# Row 82
$ws.Range("H82").Value = 3031.2593
$ws.Range("I82").Value = 1762.25
$ws.Range("J82").Value = 4877.091
$ws.Range("K82").Value = 1762.25
$ws.Range("L82").Value = 4877.091
$ws.Range("M82").Value = -1401.25
$ws.Range("N82").Value = -5599.091

# Row 85
$ws.Range("H85").Value = 3031.2593
$ws.Range("I85").Value = 1762.25
$ws.Range("J85").Value = 4877.091
$ws.Range("K85").Value = 1762.25
$ws.Range("L85").Value = 4877.091
$ws.Range("M85").Value = -514.25
$ws.Range("N85").Value = -7373.091

# Row 122
$ws.Range("H122").Value = 804012.0600000001
$ws.Range("I122").Value = 1819740.9
$ws.Range("J122").Value = 5939.357
$ws.Range("K122").Value = 5459222.699999999
$ws.Range("L122").Value = 17818.071
$ws.Range("M122").Value = -5456772.699999999
$ws.Range("N122").Value = -22718.071

# Row 126
$ws.Range("H126").Value = 4682.171
$ws.Range("I126").Value = 4731.963
$ws.Range("J126").Value = 4586.143
$ws.Range("K126").Value = 14195.889
$ws.Range("L126").Value = 13758.429
$ws.Range("M126").Value = -11725.889

# Row 132
$ws.Range("H132").Value = 2760.5
$ws.Range("I132").Value = 2458.2693
$ws.Range("J132").Value = 3742.75
$ws.Range("K132").Value = 7374.8079
$ws.Range("L132").Value = 11228.25
$ws.Range("M132").Value = -4844.8079
$ws.Range("N132").Value = -16288.25

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10371.637
$ws.Range("I62").Value = 7370.8887
$ws.Range("J62").Value = 23875
$ws.Range("K62").Value = 7370.8887
$ws.Range("L62").Value = 23875
$ws.Range("M62").Value = -6746.8887
$ws.Range("N62").Value = -25123

# Row 65
$ws.Range("H65").Value = 10371.637
$ws.Range("I65").Value = 7370.8887
$ws.Range("J65").Value = 23875
$ws.Range("K65").Value = 36854.4435
$ws.Range("L65").Value = 119375
$ws.Range("M65").Value = -33734.4435
$ws.Range("N65").Value = -125615

# Row 96
$ws.Range("H96").Value = 1584.5385
$ws.Range("I96").Value = 1365.9166
$ws.Range("J96").Value = 1771.9286
$ws.Range("K96").Value = 1365.9166
$ws.Range("L96").Value = 1771.9286
$ws.Range("M96").Value = 7.083399999999983
